# Add a new "abilities" column (K) to the hero table so abilities can be
# applied to heroes by table (see commit message "apply ability to hero by
# table").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the data values first - put the row4/row5 ability key in before the
# header cells so the new shared-string table ends up in the same order as
# authored (Level:...LaserStrike first -> index 33, then the three header
# strings -> 34/35/36).
$ws.Range("K4").Value = "Level:ActorConfigs:AbilityConfigHeroLaserStrike"
$ws.Range("K5").Value = "Level:ActorConfigs:AbilityConfigHeroLaserStrike"

$ws.Range("K1").Value = "abilities"
$ws.Range("K2").Value = "(array#sep=,),string"
$ws.Range("K3").Value = "ability config key"

# Match the formatting of the adjacent "prefab" column (J): header rows use
# the existing shaded styles, and the thick-bottom-border style on row 3.
$ws.Range("J1:J3").Copy() | Out-Null
$ws.Range("K1:K3").PasteSpecial(-4122) | Out-Null

# Size the new column similarly to the other wide text columns.
$ws.Columns.Item(11).ColumnWidth = 42.714286

# Restore the active selection to match where the author ended up after
# editing the new column.
$ws.Range("K13").Select() | Out-Null
